$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Rows 28-29: switch Approved -> Rejected, add ReasonToReject "Nil"
$ws.Range("I28").Value = "Rejected"
$ws.Range("J28").Value = "Nil"
$ws.Range("I29").Value = "Rejected"
$ws.Range("J29").Value = "Nil"

# Rows 30-32: switch Rejected/Nil -> Approved, clear ReasonToReject
$ws.Range("I30").Value = "Approved"
$ws.Range("J30").ClearContents()
$ws.Range("I31").Value = "Approved"
$ws.Range("J31").ClearContents()
$ws.Range("I32").Value = "Approved"
$ws.Range("J32").ClearContents()

# Update the selection to match the new active cell/selection
$ws.Range("I32:J32").Select()
